$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.587.30'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.918.73'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.47'
$ws.Range("E5").Value = '  +3.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06839'
$ws.Range("E9").Value = '  +4.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '104.59'
$ws.Range("E10").Value = '  -5.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.30'
$ws.Range("E11").Value = '  -3.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.914.96'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07690'
$ws.Range("E13").Value = '  +2.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.266'
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6688'
$ws.Range("E15").Value = '  +4.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '288.89'
$ws.Range("E16").Value = '  -10.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.585.74'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007589'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.91'
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.167.71'
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.498'
$ws.Range("E22").Value = '  +7.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.296'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.17'
$ws.Range("E26").Value = '  +1.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.99'
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.118'
$ws.Range("E28").Value = '  +6.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1066'
$ws.Range("E29").Value = '  -2.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.391'
$ws.Range("E30").Value = '  +4.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.167'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.075'
$ws.Range("E32").Value = '  +3.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05030'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7346'
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.146'
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02058'
$ws.Range("E36").Value = '  +6.02%  '
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.055'
$ws.Range("E39").Value = '  +2.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.24'
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8790'
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4393'
$ws.Range("E42").Value = '  +6.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.873'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.96'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.255'
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.309'
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.66'
$ws.Range("E48").Value = '  +11.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1229'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.81'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4041'
$ws.Range("E51").Value = '  +6.61%  '
